# Updates cryptos list values (prices, names/links, volume %) per source commit diff.
# Note: some "Price" column values look numeric (e.g. "303.17"); Excel COM auto-converts
# such strings to numbers unless the cell is pre-formatted as Text ("@"), so we set
# NumberFormat = "@" first for those cells to keep them as text, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.435.63"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "2.328.34"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.17"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.21"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.60"
$ws.Range("E11").Value = "  +7.52%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "2.688.65"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "2.317.60"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "43.298.32"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.09"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.97"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.51"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E24").Value = "  +4.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.10"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.52"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.48"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.03"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.90"
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("E35").Value = "  -7.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0704"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.80"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "1.992.58"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.74"
$ws.Range("E43").Value = "  +6.47%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.20"
$ws.Range("E45").Value = "  +3.51%  "
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.90"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.06"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.555.26"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.76"
$ws.Range("E51").Value = "  +0.68%  "
